# Append a new run-log row to the Nalco PDF tracking sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the next blank row just after the current used range (header row 1
# plus any existing log rows), so this keeps working as the log grows.
$lastRow = $ws.UsedRange.Rows.Count
$row = $lastRow + 1

$ws.Cells.Item($row, 1).Value = "2025-08-13 06:51:52 UTC"
$ws.Cells.Item($row, 2).Value = "2025-08-13 12:21:52 IST"
$ws.Cells.Item($row, 3).Value = "SKIPPED"
$ws.Cells.Item($row, 4).Value = "No change in PDF. Skipping download & Excel update."
$ws.Cells.Item($row, 5).Value = "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-07-08-2025.pdf"
# Column F (Saved PDF) and H (Total Rows After) stay blank - nothing was
# downloaded/appended this run.
$ws.Cells.Item($row, 7).Value = 0

# Match the centered alignment style used by every other data row.
$newRowRange = $ws.Range($ws.Cells.Item($row, 1), $ws.Cells.Item($row, 8))
$newRowRange.HorizontalAlignment = -4108 # xlCenter
$newRowRange.VerticalAlignment = -4108 # xlCenter
